# Weekly update: a new price-report row is inserted at row 9 (for
# "Vega Monumental Concepción" / Haba), pushing the existing rows 9-33
# down to 10-34. The rest of the data is unchanged, just shifted down
# by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 9, shifting rows 9:33 down to 10:34.
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 44811
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 100112026
$ws.Range("G9").Value = "Haba"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 12600
$ws.Range("N9").Value = "`$/saco 25 kilos"
$ws.Range("O9").Value = "Región de Coquimbo"
$ws.Range("P9").Value = 504
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
